$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (Changed) date, stored as an Excel serial
# date. Every populated data row (2 through 536) had its value bumped by
# exactly one day (45225 -> 45226).
for ($r = 2; $r -le 536; $r++) {
    $ws.Cells.Item($r, 3).Value = 45226
}
